$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Header row, re-using the same bold/bordered/centered style as the other sheets ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$wsForecast.Range("A1:B1").Copy()
$wsForecast.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- 4. Data rows (83 rows -> worksheet rows 2..84) ---
$rows = @(
    @(44948.99999999999,229,-69.39903885156741,531.6907136245158),
    @(44955.99999999999,230,-61.6455941943687,518.3512123313933),
    @(44976.99999999999,234,-63.19829992924354,547.2389613622099),
    @(44983.99999999999,236,-71.95501298682674,554.8635409760451),
    @(44990.99999999999,237,-73.26182167267376,530.2531739143109),
    @(44997.99999999999,239,-60.63540269055888,536.1836159734152),
    @(45004.99999999999,240,-51.18540547179876,534.1221911923178),
    @(45011.99999999999,241,-42.7957287817549,537.9868847422954),
    @(45039.99999999999,247,-47.63151175171982,538.336363498788),
    @(45046.99999999999,248,-49.48214608043078,559.0073594194529),
    @(45053.99999999999,249,-34.17044742102082,543.4528681359083),
    @(45060.99999999999,251,-59.07822566571439,564.4654752657818),
    @(45067.99999999999,252,-48.68203664301404,525.4352581699796),
    @(45074.99999999999,253,-51.37646101153899,551.1131341580498),
    @(45081.99999999999,255,-60.28799659025523,544.6683959039249),
    @(45088.99999999999,256,-36.82324723976456,551.1389206936187),
    @(45095.99999999999,257,-38.35250485660982,576.2813267552229),
    @(45102.99999999999,259,-41.54495436883936,523.5489449966459),
    @(45109.99999999999,260,-42.43212503869997,557.5233348678169),
    @(45130.99999999999,264,-28.42970129903944,577.1798130245559),
    @(45144.99999999999,267,-42.66175958813193,561.5569087723908),
    @(45151.99999999999,268,-43.88663082095341,550.6659942136257),
    @(45165.99999999999,271,-12.66700960368427,569.1275699797906),
    @(45172.99999999999,272,-19.31540997581291,568.1151138508686),
    @(45179.99999999999,274,-24.61378605058202,566.1832398913858),
    @(45186.99999999999,275,-15.93655626184784,586.9451901181714),
    @(45193.99999999999,276,-25.61956062731792,568.8678490458906),
    @(45200.99999999999,278,-45.34786025272531,574.0940649117211),
    @(45207.99999999999,279,-44.85655546870279,604.6772114418445),
    @(45214.99999999999,280,-38.93630620728349,591.4680131006191),
    @(45221.99999999999,282,-22.42416199014703,591.9129174589186),
    @(45228.99999999999,283,-20.76263727128046,584.3643455245571),
    @(45235.99999999999,284,-17.76144537256967,587.7284504303012),
    @(45242.99999999999,286,-25.90867400347904,582.0484203054766),
    @(45249.99999999999,287,-31.88460263248604,602.0707161712329),
    @(45256.99999999999,288,-22.72656956475656,590.9483230770186),
    @(45263.99999999999,290,-15.76604666346423,572.2957506821484),
    @(45270.99999999999,291,-2.43667147890208,582.7579331919295),
    @(45277.99999999999,293,-5.209959788127661,601.2474128639579),
    @(45298.99999999999,297,-7.484876477867366,603.3930782648101),
    @(45312.99999999999,299,11.88663564351645,595.4807833177304),
    @(45326.99999999999,302,1.502534538708446,611.799030944981),
    @(45333.99999999999,303,1.5376849573836,626.0626020449951),
    @(45347.99999999999,306,25.96229968826541,609.6968622581526),
    @(45361.99999999999,309,23.28608760466798,624.6275446412997),
    @(45368.99999999999,310,14.21026573541645,594.9225754802721),
    @(45375.99999999999,311,11.55589515604226,614.0985869089287),
    @(45382.99999999999,313,6.663252099787602,622.7491042729168),
    @(45389.99999999999,314,6.423058311196197,621.0876492794581),
    @(45396.99999999999,316,12.04217557829045,614.3330817757487),
    @(45403.99999999999,317,11.48237577995817,633.6706654991186),
    @(45410.99999999999,318,22.79538841408689,620.095236850163),
    @(45417.99999999999,320,0.9530848273931074,638.2975476903929),
    @(45424.99999999999,321,28.13885631303251,640.3436769571459),
    @(45431.99999999999,322,16.28739247320754,609.5700361908247),
    @(45445.99999999999,325,4.616037038658737,646.3277165757827),
    @(45452.99999999999,326,30.70104794762134,642.6329572224414),
    @(45459.99999999999,328,28.86542775127044,629.6879650089012),
    @(45466.99999999999,329,39.52205969015039,640.0224570781299),
    @(45473.99999999999,330,30.01114870428676,638.4647543895766),
    @(45480.99999999999,332,21.38487963003394,647.3637745966629),
    @(45487.99999999999,333,17.17673006626303,629.6004057702796),
    @(45501.99999999999,336,30.73888686166037,620.9181941757207),
    @(45508.99999999999,337,23.64170315504066,632.8894137468856),
    @(45515.99999999999,338,32.15970155383832,635.4702440135416),
    @(45529.99999999999,341,42.55938242666809,657.997676675741),
    @(45536.99999999999,343,44.66768549224246,638.8623572775244),
    @(45543.99999999999,344,50.24115631039649,648.0068961073488),
    @(45550.99999999999,345,36.31946180785185,645.4167626891121),
    @(45557.99999999999,347,37.06406814674168,654.9731021490725),
    @(45564.99999999999,348,10.81331906207743,670.4642238696215),
    @(45571.99999999999,349,43.61625880571886,641.4685567584),
    @(45578.99999999999,351,46.25845969486508,644.4290188132682),
    @(45585.99999999999,352,59.68560036714276,640.3313271075068),
    @(45599.99999999999,355,44.61301117582777,661.5939652461927),
    @(45606.99999999999,356,61.58820094464421,681.3829829633728),
    @(45613.99999999999,357,73.3919256160103,658.45337887372),
    @(45620.99999999999,359,40.43513812674061,662.961485305921),
    @(45627.99999999999,360,55.30912650205082,655.6277277847822),
    @(45634.99999999999,361,72.04472097301804,663.2391900905236),
    @(45641.99999999999,363,65.79188384707582,701.3819250975336),
    @(45648.99999999999,364,63.1633785833931,645.6708884888151),
    @(45655.99999999999,366,43.93508478977477,668.7600443474503)
)

$data = New-Object 'object[,]' $rows.Count, 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$wsForecast.Range("A2:D84").Value = $data

# --- 5. Re-apply the date style (column A) used on the other sheets ---
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A84").PasteSpecial(-4122)  # xlPasteFormats

# Re-write the data again, since PasteSpecial(formats) should not disturb values,
# but make sure values are intact after style paste.
$wsForecast.Range("A2:D84").Value = $data

Write-Host "PO Forecast sheet created with $($rows.Count) data rows."
